$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'237.71"
$ws.Range("D3").Value = "'21.76"
$ws.Range("D4").Value = "'5.386"
$ws.Range("D5").Value = "'0.05625"
$ws.Range("D6").Value = "'6.478"
$ws.Range("D7").Value = "'3.350"
$ws.Range("D8").Value = "'0.7967"
$ws.Range("D9").Value = "'1.030"
$ws.Range("D10").Value = "'0.1387"
$ws.Range("D11").Value = "'0.07307"
$ws.Range("D12").Value = "'0.03114"
$ws.Range("D13").Value = "'0.02966"
$ws.Range("D14").Value = "'0.09219"
$ws.Range("D15").Value = "'0.001663"
$ws.Range("D16").Value = "'3.256"
$ws.Range("D17").Value = "'0.04775"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006231"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.005074"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001051"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "UpBots"
$ws.Range("C22").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D22").Value = "'0.0003901"
$ws.Range("E22").Value = "21UpBotsUBXT"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.932"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.201"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D25").Value = "'0.01167"
$ws.Range("E25").Value = "24OneONEBestin24h"
$ws.Range("D40").Value = "'0.04077"
$ws.Range("D41").Value = "'0.006947"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003501"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1039"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.008825"
$ws.Range("D45").Value = "'0.00005435"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.6754"
$ws.Range("D48").Value = "'0.03720"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.01010"
